# Apply the latest cryptos.xlsx price/volume refresh (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell and its new literal value. Numeric-looking
# Price values are apostrophe-prefixed so Excel stores them as literal text
# (matching the columns existing text format) instead of silently coercing
# them to numbers and dropping significant trailing zeros, e.g. "17.80" -> 17.8.
$updates = @(
    @{ Cell = 'D2'; Value = '68.228.65' }
    @{ Cell = 'E2'; Value = '  -0.83%  ' }
    @{ Cell = 'D3'; Value = '2.645.97' }
    @{ Cell = 'E3'; Value = '  -0.52%  ' }
    @{ Cell = 'E4'; Value = '  +0.01%  ' }
    @{ Cell = 'D5'; Value = '''597.61' }
    @{ Cell = 'E5'; Value = '  -0.50%  ' }
    @{ Cell = 'D6'; Value = '''156.65' }
    @{ Cell = 'E6'; Value = '  +0.63%  ' }
    @{ Cell = 'E7'; Value = '  +0.03%  ' }
    @{ Cell = 'E8'; Value = '  -0.85%  ' }
    @{ Cell = 'E9'; Value = '  +1.64%  ' }
    @{ Cell = 'E10'; Value = '  -1.23%  ' }
    @{ Cell = 'E11'; Value = '  +0.26%  ' }
    @{ Cell = 'E12'; Value = '  +0.42%  ' }
    @{ Cell = 'D13'; Value = '''28.02' }
    @{ Cell = 'E13'; Value = '  -0.20%  ' }
    @{ Cell = 'E14'; Value = '  +0.73%  ' }
    @{ Cell = 'D15'; Value = '3.126.74' }
    @{ Cell = 'E15'; Value = '  -0.47%  ' }
    @{ Cell = 'D16'; Value = '68.229.41' }
    @{ Cell = 'E16'; Value = '  -0.68%  ' }
    @{ Cell = 'D17'; Value = '2.644.67' }
    @{ Cell = 'E17'; Value = '  -0.51%  ' }
    @{ Cell = 'D18'; Value = '''11.39' }
    @{ Cell = 'E18'; Value = '  -0.54%  ' }
    @{ Cell = 'D19'; Value = '''363.92' }
    @{ Cell = 'E19'; Value = '  -0.91%  ' }
    @{ Cell = 'D20'; Value = '''7.34' }
    @{ Cell = 'E20'; Value = '  -1.49%  ' }
    @{ Cell = 'E21'; Value = '  +3.17%  ' }
    @{ Cell = 'D22'; Value = '''4.79' }
    @{ Cell = 'E22'; Value = '  -1.88%  ' }
    @{ Cell = 'D23'; Value = '''2.08' }
    @{ Cell = 'E23'; Value = '  -3.06%  ' }
    @{ Cell = 'D24'; Value = '''75.52' }
    @{ Cell = 'E24'; Value = '  +3.90%  ' }
    @{ Cell = 'E25'; Value = '  -0.04%  ' }
    @{ Cell = 'E26'; Value = '  -2.60%  ' }
    @{ Cell = 'B27'; Value = 'WrappedeETH' }
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth' }
    @{ Cell = 'D27'; Value = '2.776.57' }
    @{ Cell = 'E27'; Value = '  -0.41%  ' }
    @{ Cell = 'B28'; Value = 'Binance-PegBSC-USD' }
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd' }
    @{ Cell = 'D28'; Value = '''1.03' }
    @{ Cell = 'E28'; Value = '  +3.10%  ' }
    @{ Cell = 'E29'; Value = '  -1.47%  ' }
    @{ Cell = 'D30'; Value = '''556.14' }
    @{ Cell = 'E30'; Value = '  -3.68%  ' }
    @{ Cell = 'D31'; Value = '''8.06' }
    @{ Cell = 'E31'; Value = '  +0.81%  ' }
    @{ Cell = 'E32'; Value = '  -0.76%  ' }
    @{ Cell = 'E33'; Value = '  -0.76%  ' }
    @{ Cell = 'E34'; Value = '  +0.01%  ' }
    @{ Cell = 'E35'; Value = '  -2.09%  ' }
    @{ Cell = 'D36'; Value = '''1.55' }
    @{ Cell = 'E36'; Value = '  +0.05%  ' }
    @{ Cell = 'D37'; Value = '''161.54' }
    @{ Cell = 'E37'; Value = '  +1.38%  ' }
    @{ Cell = 'D38'; Value = '''19.83' }
    @{ Cell = 'E38'; Value = '  +2.74%  ' }
    @{ Cell = 'E39'; Value = '  +0.92%  ' }
    @{ Cell = 'E40'; Value = '  -3.38%  ' }
    @{ Cell = 'D41'; Value = '''5.33' }
    @{ Cell = 'E41'; Value = '  -1.66%  ' }
    @{ Cell = 'B42'; Value = 'BabyDogeCoin' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge' }
    @{ Cell = 'D42'; Value = '0.0₆0334' }
    @{ Cell = 'E42'; Value = '  +3.26%  ' }
    @{ Cell = 'B43'; Value = 'WhiteBITCoin' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt' }
    @{ Cell = 'D43'; Value = '''17.80' }
    @{ Cell = 'E43'; Value = '  +0.32%  ' }
    @{ Cell = 'E44'; Value = '  -1.93%  ' }
    @{ Cell = 'D46'; Value = '''158.80' }
    @{ Cell = 'E46'; Value = '  +1.30%  ' }
    @{ Cell = 'D47'; Value = '''3.73' }
    @{ Cell = 'E47'; Value = '  -0.38%  ' }
    @{ Cell = 'E48'; Value = '  -0.09%  ' }
    @{ Cell = 'B49'; Value = 'Optimism' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op' }
    @{ Cell = 'D49'; Value = '''1.69' }
    @{ Cell = 'E49'; Value = '  -2.08%  ' }
    @{ Cell = 'B50'; Value = 'Cronos' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' }
    @{ Cell = 'D50'; Value = '''0.0783' }
    @{ Cell = 'E50'; Value = '  +0.28%  ' }
    @{ Cell = 'D51'; Value = '''0.615' }
    @{ Cell = 'E51'; Value = '  -0.65%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
